$p = $ppt.ActivePresentation

# The three "CropOutPadding" demo slides each have a small "Expected
# Output" caption text box that was drawn *behind* the picture(s) it
# labels. A newly-added test now checks shape z-order on this deck, so
# bring the caption text box to the front (== push the picture(s) to
# the back) on each of the affected slides.

# Slide 5: "text 3" (caption) + "selectMe" (picture)
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item("text 3").ZOrder(0)   # msoBringToFront

# Slide 5 also picks up a slow, 2-second slide transition.
$t5 = $s5.SlideShowTransition
$t5.Duration = 2
$t5.Speed = 1                          # ppTransitionSpeedSlow

# Slide 8: "text 3" (caption) + "selectMe2" + "selectMe1" (pictures)
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item("text 3").ZOrder(0)   # msoBringToFront

# Slide 11: "text 3" (caption) + "selectMe" (picture)
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item("text 3").ZOrder(0)  # msoBringToFront
